# Generate Report for Handoff
# A new handoff xliff was generated for the "8e0efb33-..." source file, so its
# "Latest Handoff Datetime" / "Latest HO Xliff Generate Date" timestamps are
# refreshed on the Overview sheet and on each per-locale handoff report sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4").Value = "2016-10-18 12:04:30"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H4").Value = "2016-10-18 12:04:19"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H4").Value = "2016-10-18 12:04:30"
